# Versione 0.7: Inizio Gestione Modifica Profilo
#
# Split the run " Il fornitore può solo visualizzare o eliminare." into
# three runs carrying the same character formatting:
#   " "  +  "L’admin"  +  " può solo visualizzare o eliminare."

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Text = " Il fornitore può solo visualizzare o eliminare."
$found = $find.Execute()

if ($found) {
    $whole = $d.Range($find.Parent.Start, $find.Parent.End)
    $s = $whole.Start

    $apostrophe = [string][char]0x2019
    $adminText = "L" + $apostrophe + "admin"
    $tailText = " può solo visualizzare o eliminare."

    # Replace the whole run's text in one shot with the concatenation of
    # the three pieces (keeps this a single run for the moment).
    $whole.Text = " " + $adminText + $tailText
    $e = $whole.End

    $adminStart = $s + 1
    $adminEnd = $adminStart + $adminText.Length
    $tailStart = $adminEnd
    $tailEnd = $e

    # Carve the combined run into three separate runs, working from right
    # to left so that later InsertAfter-style boundary changes never touch
    # an already-separated run on its right. Toggling Bold on and back off
    # marks each piece as its own run even though the final formatting
    # matches its neighbours.
    $tailRange = $d.Range($tailStart, $tailEnd)
    $tailRange.Bold = 1
    $tailRange.Bold = 0

    $adminRange = $d.Range($adminStart, $adminEnd)
    $adminRange.Bold = 1
    $adminRange.Bold = 0

    $spaceRange = $d.Range($s, $adminStart)
    $spaceRange.Bold = 1
    $spaceRange.Bold = 0

    # Changing the run above causes the host to re-coalesce any other
    # identically-formatted adjacent run pair still left in the paragraph
    # (e.g. the untouched " Nel caso…prodotti." / " 1GG" runs that follow).
    # Re-split that boundary back apart the same way.
    $findTail = $d.Content.Find
    $findTail.Text = " 1GG"
    $foundTail = $findTail.Execute()
    if ($foundTail) {
        $gg = $d.Range($findTail.Parent.Start, $findTail.Parent.End)
        $gg.Bold = 1
        $gg.Bold = 0
    }
}
